# Insert a new, blank slide as slide 2 (right after the existing slide).
#
# This mirrors the target OOXML diff:
#   - presentation.xml: sldIdLst gains a new <p:sldId id="257" .../> entry
#   - a new, essentially empty ppt/slides/slideN.xml part is added
#     (just the default empty group shape / spTree, no placeholders or
#     other shapes left on it)

$p = $ppt.ActivePresentation

# Find the "Blank" layout on the slide master (falls back to the
# well-known PpSlideLayout enum value ppLayoutBlank = 12 if, for some
# reason, no layout is named "Blank").
$blankLayout = $null
$layouts = $p.SlideMaster.CustomLayouts
for ($i = 1; $i -le $layouts.Count; $i++) {
    $candidate = $layouts.Item($i)
    if ($candidate.Name -eq "Blank") {
        $blankLayout = $candidate
        break
    }
}

# Insert the new slide immediately after the current last slide.
$newIndex = $p.Slides.Count + 1

if ($blankLayout -ne $null) {
    $newSlide = $p.Slides.AddSlide($newIndex, $blankLayout)
} else {
    $newSlide = $p.Slides.Add($newIndex, 12)
}
